# Add data for 2021-09-26
# Updates the "through" date in the sheet name + header label, and bumps
# the carjacking counts for the neighborhoods/months affected by the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab (Through 2021-09-17 -> Through 2021-09-18)
$ws.Name = "Through 2021-09-18"

# Update the column header label that mirrors the sheet name
$ws.Range("B1").Value = "September 2021 (through September 18)"

# --- Column B (September 2021, current month) increments ---
$ws.Range("B2").Value = 11   # Garfield Park
$ws.Range("B3").Value = 5    # North Lawndale
$ws.Range("B5").Value = 8    # Austin
$ws.Range("B10").Value = 3   # West Town
$ws.Range("B23").Value = 2   # United Center
$ws.Range("B27").Value = 3   # Avalon Park
$ws.Range("B32").Value = 2   # Chicago Lawn

# --- Other month/neighborhood cell updates ---
$ws.Range("BD2").Value = 2   # Garfield Park, September 2015
$ws.Range("AC3").Value = 2   # North Lawndale, September 2018
$ws.Range("K5").Value = 7    # Austin, September 2020
$ws.Range("BD13").Value = 1  # Chatham, September 2015 (new)
$ws.Range("T19").Value = 2   # Wicker Park, September 2019
$ws.Range("AC20").Value = 2  # Englewood, September 2018
$ws.Range("BD24").Value = 1  # Ashburn, September 2015 (new)
$ws.Range("AU28").Value = 1  # Logan Square, September 2016 (new)
$ws.Range("AU37").Value = 1  # Gage Park, September 2016 (new)
$ws.Range("T40").Value = 1   # West Elsdon, September 2019 (new)
$ws.Range("K55").Value = 5   # Grand Crossing, September 2020
$ws.Range("AL64").Value = 1  # Brighton Park, September 2017 (new)
$ws.Range("AU99").Value = 1  # West Ridge, September 2016 (new)
